# Weekly update: a new week of data (2023-12-12) is added at the top of the
# Ciboulette price table (rows 16-17), pushing all the existing data rows
# down by two rows (old 16..49 become 18..51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 16, shifting everything
# below (old rows 16-49) down to rows 18-51.
$ws.Rows("16:17").Insert()

# Row 16: new "Primera" quality record for 2023-12-12 (serial 45272)
$ws.Cells.Item(16, 1).Value = 7
$ws.Cells.Item(16, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(16, 3).Value = "Ñuble"
$ws.Cells.Item(16, 4).Value = 45272
$ws.Cells.Item(16, 5).Value = 16
$ws.Cells.Item(16, 6).Value = 100112039
$ws.Cells.Item(16, 7).Value = "Ciboulette"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 150
$ws.Cells.Item(16, 11).Value = 2000
$ws.Cells.Item(16, 12).Value = 2000
$ws.Cells.Item(16, 13).Value = 2000
$ws.Cells.Item(16, 14).Value = "`$/docena de atados"
$ws.Cells.Item(16, 15).Value = "Región Metropolitana"
$ws.Cells.Item(16, 16).Value = 667
$ws.Cells.Item(16, 17).Value = 3
$ws.Cells.Item(16, 18).Value = "Hortaliza"

# Row 17: new "Segunda" quality record for 2023-12-12 (serial 45272)
$ws.Cells.Item(17, 1).Value = 7
$ws.Cells.Item(17, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(17, 3).Value = "Ñuble"
$ws.Cells.Item(17, 4).Value = 45272
$ws.Cells.Item(17, 5).Value = 16
$ws.Cells.Item(17, 6).Value = 100112039
$ws.Cells.Item(17, 7).Value = "Ciboulette"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Segunda"
$ws.Cells.Item(17, 10).Value = 150
$ws.Cells.Item(17, 11).Value = 1500
$ws.Cells.Item(17, 12).Value = 1500
$ws.Cells.Item(17, 13).Value = 1500
$ws.Cells.Item(17, 14).Value = "`$/docena de atados"
$ws.Cells.Item(17, 15).Value = "Región Metropolitana"
$ws.Cells.Item(17, 16).Value = 500
$ws.Cells.Item(17, 17).Value = 3
$ws.Cells.Item(17, 18).Value = "Hortaliza"
